$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ponds")
$r = $ws1.Range("K3")
Write-Host "Interior.ThemeColor: $($r.Interior.ThemeColor())"
Write-Host "Interior.TintAndShade: $($r.Interior.TintAndShade())"
Write-Host "Interior.Color: $($r.Interior.Color())"
Write-Host "Interior.Pattern: $($r.Interior.Pattern())"
Write-Host "Font.Color: $($r.Font.Color())"
Write-Host "Font.ThemeColor: $($r.Font.ThemeColor())"
Write-Host "Font.Size: $($r.Font.Size())"
Write-Host "Borders(9).Weight: $($r.Borders.Item(9).Weight())"
Write-Host "Borders(9).Color: $($r.Borders.Item(9).Color())"
Write-Host "Borders(9).ColorIndex: $($r.Borders.Item(9).ColorIndex())"
